$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "Testing" worksheet as the last (7th) tab.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$testing = $wb.Worksheets.Add($null, $lastSheet)
$testing.Name = "Testing"

# Header row (row 3)
$testing.Range("B3").Value = "Cases/million"
$testing.Range("C3").Value = "Test/million"
$testing.Range("D3").Value = "Ratio test/case"

# Data rows
$testing.Range("A4").Value = "Spain"
$testing.Range("B4").Value = 6366
$testing.Range("C4").Value = 116543

$testing.Range("A5").Value = "Germany"
$testing.Range("B5").Value = 2359
$testing.Range("C5").Value = 70100

$testing.Range("A6").Value = "Sweden"
$testing.Range("B6").Value = 7071
$testing.Range("C6").Value = 51397

$testing.Range("A7").Value = "Ireland"
$testing.Range("B7").Value = 5169
$testing.Range("C7").Value = 94738

$testing.Range("A8").Value = "England and Wales (UK)"
$testing.Range("B8").Value = 4209
$testing.Range("C8").Value = 156889

# Ratio formulas: D4 stands alone, D5:D8 becomes a shared-formula block
$testing.Range("D4").Formula = "=C4/B4"
$testing.Range("D5:D8").Formula = "=C5/B5"

# 2 decimal place number format for the ratio column
$testing.Range("D4:D8").NumberFormat = "0.00"

# Column width + page setup to match the authored sheet
$testing.Columns.Item(4).ColumnWidth = 13.92
$testing.PageSetup.PaperSize = 9
$testing.PageSetup.Orientation = 1

# View state: D8 is the active cell of a D3:D8 selection
$testing.Range("D3:D8").Select()

# ------------------------------------------------------------------
# 2. UK sheet: add a "10th death" column before the existing
#    "100th death" column (shifts that single header cell from D1 to E1).
# ------------------------------------------------------------------
$uk = $wb.Worksheets.Item("UK")
$uk.Range("D1").Copy($uk.Range("E1"))
$uk.Range("D1").Value = "10th death"

$uk.Columns.Item(2).ColumnWidth = 16.92
$uk.Columns.Item(3).ColumnWidth = 19.76

# Restore UK as the active/selected sheet with the recorded selection.
$uk.Activate()
$uk.Range("C6").Select()
